$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.23336480490412725
$ws.Range("A2").Value = -0.022905092534752924
$ws.Range("A3").Value = -0.0039999999866324742
$ws.Range("A4").Value = -0.0079999999754836182
$ws.Range("A5").Value = -0.0029999999861587412
$ws.Range("A6").Value = -0.0019999999851023631
$ws.Range("A7").Value = -0.0099949464808917199
$ws.Range("A8").Value = -0.0099999999640472126
$ws.Range("A9").Value = -0.0019999999819271252
$ws.Range("A10").Value = -0.001999999980537126
$ws.Range("A11").Value = -0.0029999999779386499
$ws.Range("A12").Value = -0.003499999976356083
$ws.Range("A13").Value = -0.0034999999749363297
$ws.Range("A14").Value = -0.0079999999636175545
$ws.Range("A15").Value = -0.00099999998029431936
$ws.Range("A16").Value = -0.0019999999776141308
$ws.Range("A17").Value = -0.0019999999772455368
$ws.Range("A18").Value = -0.0039999999722892809
$ws.Range("A19").Value = -0.0039999999893565175
$ws.Range("A20").Value = -0.0039999999861741742
$ws.Range("A21").Value = -0.0039999999854520851
$ws.Range("A22").Value = -0.040479060767984798
$ws.Range("A23").Value = -0.0049999999834327014
$ws.Range("A24").Value = -0.019999999945747859
$ws.Range("A25").Value = -0.019999999945055968
$ws.Range("A26").Value = -0.0024999999836268216
$ws.Range("A27").Value = -0.0024999999833075215
$ws.Range("A28").Value = -0.0019999999830648818
$ws.Range("A29").Value = -0.0069999999699437154
$ws.Range("A30").Value = -0.0030681049983933839
$ws.Range("A31").Value = -0.0069999999686292114
$ws.Range("A32").Value = -0.0099999999612236934
$ws.Range("A33").Value = -0.0039999999756759053
